# "Generate Report for Handoff" - b.md went through a new handoff cycle:
# a fresh handoff xlf was generated for both zh-cn and de-de locales, so
# the status flips from "Handed back: in sync with en-US" to
# "Ready for handoff", new handoff files/timestamps are recorded, and an
# error is surfaced because the existing handback isn't for the latest
# source version yet.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3061528868bf82d47e4311f1f8f2ab9c52a62dbf/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/080524024c4031ee2e718accca374fe90e4cb051/e2e/b.md."

# --- Overview sheet: row 3 is b.md -------------------------------------
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E3").Value = "Ready for handoff"
$ov.Range("F3").Value = "Ready for handoff"
$ov.Range("G3").Value = "2016-08-18 12:37:18"

# --- zh-cn sheet: row 3 is b.md ------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("F3").Value = "False"
$zh.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zh.Range("H3").Value = "2016-08-18 12:37:13"
$zh.Range("P3").Value = $errorDetail
$zh.Columns.Item(16).ColumnWidth = 39.166666666666664

# --- de-de sheet: row 3 is b.md ------------------------------------------
$de = $wb.Worksheets.Item("de-de")
$de.Range("C3").Value = "Ready for handoff"
$de.Range("F3").Value = "False"
$de.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$de.Range("H3").Value = "2016-08-18 12:37:18"
$de.Range("P3").Value = $errorDetail
$de.Columns.Item(16).ColumnWidth = 39.166666666666664
